# Update "想去人数" (interest count) figures in column F across the
# three affected worksheets: 展览, 本地生活, 全部类型.
# (演出 sheet is untouched by this update.)

$wb = $excel.ActiveWorkbook

# --- 展览 (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 61
$ws1.Range("F3").Value = 1171
$ws1.Range("F5").Value = 70
$ws1.Range("F7").Value = 830
$ws1.Range("F8").Value = 434
$ws1.Range("F10").Value = 2083
$ws1.Range("F12").Value = 252
$ws1.Range("F14").Value = 971
$ws1.Range("F16").Value = 2105
$ws1.Range("F17").Value = 575
$ws1.Range("F18").Value = 10186
$ws1.Range("F19").Value = 993
$ws1.Range("F20").Value = 540
$ws1.Range("F21").Value = 99
$ws1.Range("F22").Value = 122
$ws1.Range("F24").Value = 250

# --- 本地生活 (Local life) ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 5651
$ws3.Range("F3").Value = 454
$ws3.Range("F4").Value = 430

# --- 全部类型 (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 61
$ws4.Range("F3").Value = 5651
$ws4.Range("F4").Value = 454
$ws4.Range("F5").Value = 430
$ws4.Range("F7").Value = 1171
$ws4.Range("F10").Value = 70
$ws4.Range("F12").Value = 830
$ws4.Range("F14").Value = 434
$ws4.Range("F16").Value = 2083
$ws4.Range("F18").Value = 252
$ws4.Range("F22").Value = 971
$ws4.Range("F27").Value = 2105
$ws4.Range("F28").Value = 575
$ws4.Range("F31").Value = 993
$ws4.Range("F32").Value = 540
$ws4.Range("F33").Value = 99
$ws4.Range("F34").Value = 122
$ws4.Range("F39").Value = 250
